# Updated cryptos list on Mon May 27 11:37:25 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "68.510.14"
Set-TextValue $ws.Range("E2") "  -0.89%  "
Set-TextValue $ws.Range("D3") "3.901.10"
Set-TextValue $ws.Range("E3") "  +2.33%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "602.93"
Set-TextValue $ws.Range("E5") "  +0.18%  "
Set-TextValue $ws.Range("D6") "167.58"
Set-TextValue $ws.Range("E6") "  +2.15%  "
Set-TextValue $ws.Range("D7") "3.900.16"
Set-TextValue $ws.Range("E7") "  +2.34%  "
Set-TextValue $ws.Range("E8") "  +0.32%  "
Set-TextValue $ws.Range("E9") "  -1.26%  "
Set-TextValue $ws.Range("D11") "6.48"
Set-TextValue $ws.Range("E11") "  +2.88%  "
Set-TextValue $ws.Range("E13") "  +3.77%  "
Set-TextValue $ws.Range("D14") "37.51"
Set-TextValue $ws.Range("E14") "  +0.63%  "
Set-TextValue $ws.Range("D15") "4.556.27"
Set-TextValue $ws.Range("E15") "  +2.46%  "
Set-TextValue $ws.Range("D16") "3.885.42"
Set-TextValue $ws.Range("E16") "  +1.98%  "
Set-TextValue $ws.Range("D17") "68.633.25"
Set-TextValue $ws.Range("E17") "  -0.90%  "
Set-TextValue $ws.Range("E18") "  +0.39%  "
Set-TextValue $ws.Range("D19") "17.31"
Set-TextValue $ws.Range("E19") "  -0.17%  "
Set-TextValue $ws.Range("D20") "0.112"
Set-TextValue $ws.Range("E20") "  -2.13%  "
Set-TextValue $ws.Range("E21") "  -3.72%  "
Set-TextValue $ws.Range("D22") "489.55"
Set-TextValue $ws.Range("E22") "  +0.21%  "
Set-TextValue $ws.Range("D23") "0.726"
Set-TextValue $ws.Range("E23") "  +0.31%  "
Set-TextValue $ws.Range("E24") "  +4.07%  "
Set-TextValue $ws.Range("D25") "84.68"
Set-TextValue $ws.Range("E25") "  -0.20%  "
Set-TextValue $ws.Range("E26") "  -1.19%  "
Set-TextValue $ws.Range("D27") "12.03"
Set-TextValue $ws.Range("E27") "  -1.79%  "
Set-TextValue $ws.Range("D28") "10.18"
Set-TextValue $ws.Range("E28") "  +1.36%  "
Set-TextValue $ws.Range("E29") "  +0.01%  "
Set-TextValue $ws.Range("E30") "  -1.12%  "
Set-TextValue $ws.Range("D31") "4.054.44"
Set-TextValue $ws.Range("E31") "  +2.36%  "
Set-TextValue $ws.Range("E32") "  -0.87%  "
Set-TextValue $ws.Range("E33") "  -3.63%  "
Set-TextValue $ws.Range("D34") "31.88"
Set-TextValue $ws.Range("E34") "  -0.15%  "
Set-TextValue $ws.Range("D35") "3.858.28"
Set-TextValue $ws.Range("E35") "  +2.67%  "
Set-TextValue $ws.Range("E36") "  -0.51%  "
Set-TextValue $ws.Range("D37") "1.03"
Set-TextValue $ws.Range("E37") "  +1.15%  "
Set-TextValue $ws.Range("E38") "  +0.43%  "
Set-TextValue $ws.Range("E39") "  -1.29%  "
Set-TextValue $ws.Range("D40") "3.17"
Set-TextValue $ws.Range("E40") "  +4.36%  "
Set-TextValue $ws.Range("E41") "  +0.00%  "
Set-TextValue $ws.Range("E42") "  -1.08%  "
Set-TextValue $ws.Range("D43") "431.98"
Set-TextValue $ws.Range("E43") "  +1.27%  "
Set-TextValue $ws.Range("E44") "  -0.33%  "
Set-TextValue $ws.Range("D45") "48.22"
Set-TextValue $ws.Range("E45") "  -0.77%  "
Set-TextValue $ws.Range("D46") "8.54"
Set-TextValue $ws.Range("E46") "  +1.72%  "
Set-TextValue $ws.Range("D48") "142.76"
Set-TextValue $ws.Range("E48") "  +1.04%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue $ws.Range("D49") "0.000269"
Set-TextValue $ws.Range("E49") "  +18.20%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D50") "2.804.47"
Set-TextValue $ws.Range("E50") "  -1.07%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D51") "39.45"
Set-TextValue $ws.Range("E51") "  -0.25%  "
